# Refresh nightfall data tables (waves.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ------------------------------------------------

# Row 7
$ws.Range("J7").Value = "Choir acolytes chant, boosting nearby allies."

# Row 8
$ws.Range("F8").Value = "enemy:abyssal-howler"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8"
$ws.Range("J8").Value = "Howlers arrive after fog surge; manage sanity bleed."

# Row 9
$ws.Range("F9").Value = "enemy:null-sentinel"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "5"
$ws.Range("J9").Value = "Sentinels sweep beams while shamblers close in."

# --- New row 10 ------------------------------------------------------------

$ws.Range("A10:E10").NumberFormat = "@"
$ws.Range("G10:H10").NumberFormat = "@"

$ws.Range("A10").Value = "60"
$ws.Range("B10").Value = "08"
$ws.Range("C10").Value = "0005"
$ws.Range("D10").Value = "240"
$ws.Range("E10").Value = "40"
$ws.Range("F10").Value = "enemy:myriad-fragment"
$ws.Range("G10").Value = "16"
$ws.Range("H10").Value = "8.0"
$ws.Range("I10").Value = "swarm"
$ws.Range("J10").Value = "Fragments dash in packs forcing kite routes."

# --- New row 11 ------------------------------------------------------------

$ws.Range("A11:E11").NumberFormat = "@"
$ws.Range("G11:H11").NumberFormat = "@"

$ws.Range("A11").Value = "60"
$ws.Range("B11").Value = "08"
$ws.Range("C11").Value = "0006"
$ws.Range("D11").Value = "270"
$ws.Range("E11").Value = "50"
$ws.Range("F11").Value = "enemy:harbor-dredger"
$ws.Range("G11").Value = "6"
$ws.Range("H11").Value = "10.5"
$ws.Range("I11").Value = "line"
$ws.Range("J11").Value = "Dredgers charge straight for the beacon core."

